$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column, matching the style of the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save values for each data row (2-20)
$saveValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
